$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.725.13"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.81%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.869.15"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.64%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4712"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.55%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2757"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.66%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06379"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.49%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "17.90"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +10.16%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.854.87"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.12%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07451"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.14%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.66%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "85.19"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.01%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6357"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.82%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "30.683.41"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.83%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "245.58"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +6.40%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9986"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.14%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.85"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.95%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007408"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.01%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.003"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.085"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.50%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.397"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.51%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "164.22"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.75%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.25"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +2.52%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.911"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.70%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1020"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.68%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.383"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.087"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.98%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.04964"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.158"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.56%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7093"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.71%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.711"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.82%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.01909"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.700"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +2.46%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.8852"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.46%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.010"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.88%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "105.68"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.16%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4112"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.51%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.575"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.68%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.342"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +4.28%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "65.47"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +7.28%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1224"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.39%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.707"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.06%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "33.80"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.86%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05576"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.384"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.70%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3708"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.36%  "
